$d = $word.ActiveDocument

# --- Locate the final (bookmark) paragraph; it is currently empty. ---
$n = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($n)

# Remove the _GoBack bookmark from its current spot; we'll re-add it later
# at the end of the new body paragraph.
$hadBookmark = $false
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
    $hadBookmark = $true
} catch {
    $hadBookmark = $false
}

# --- Create a brand-new paragraph *after* the (soon-to-be heading) paragraph.
# It inherits body-style paragraph-mark formatting (i val="0", sz 24) from the
# bookmark paragraph, which is exactly what we need for the new body paragraph. ---
$bookmarkPara.Range.InsertParagraphAfter()

$n2 = $d.Paragraphs.Count
$headingPara = $d.Paragraphs.Item($n2 - 1)
$bodyPara = $d.Paragraphs.Item($n2)

# --- Turn the old bookmark paragraph into the new bold heading paragraph. ---
$headingPara.Range.Text = "The State of Asynchronous Rust:"
$hr = $headingPara.Range
$hr.Font.Bold = $true
$hr.Font.Size = 16

# --- Fill in the new body paragraph with the ecosystem blurb. ---
$bodyPara.Range.Text = "The asynchronous rust ecosystem has undergone a lot of working/innovation over a period of time. The developer work hard to develop asynchronous Rust ecosystem strong. The future trait inside the standard library and the async/await language feature has recently been stabilized."

# --- Re-anchor the _GoBack bookmark to the end of the new body paragraph. ---
if ($hadBookmark) {
    $endPos = $bodyPara.Range.End
    $d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos)) | Out-Null
}

Write-Output "done"
